# faturamento_diario.xlsx update:
#  - Correct day 16 (August/2025) total_venda value.
#  - Add a new record for day 18 of August/2025, inserted right after the
#    existing day 17 row so the monthly blocks stay in order; every row
#    below shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing value for row 17 (Dia=16, Mes=8/2025)
$ws.Range("B17").Value = 21142.5

# Insert a new row at position 19 (pushes old row 19.. down to 20..)
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row with the Dia=18 / Mes=8 (08/2025) record
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 51514.1
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 2025
$ws.Range("E19").Value = "08/2025"
